$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right (B11) 6 -> 9, Wrong (C11) 3 -> 2
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 2

# Row 12 (Total): Right (B12) 108 -> 162, Wrong (C12) -12 -> -8, Max text (E12) "96/168" -> "154/252"
$ws.Range("B12").Value = 162
$ws.Range("C12").Value = -8
$ws.Range("E12").Value = "154/252"
